# "added prob 16 in common part"
# Appends a new key/buildingBlock pair (rows 40-41) to Sheet1 right after
# the existing row 39, and updates the sheet's scroll/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 40: c0030 (new "common" problem-16 key)
$ws.Range("A40").Value = "c0030"
$ws.Range("B40").Value = "로그법칙을 쓰기 위해 통일할 밑을 결정하고 밑변환공식으로 밑을 일치시킵니다."
$ws.Range("C40").Value = "밑2;"

# Row 41: c0031
$ws.Range("A41").Value = "c0031"
$ws.Range("B41").Value = "통일된 하나의 밑에 대해 로그법칙을 적용해서 로그를 계산합니다."
$ws.Range("C41").Value = "밑2;"

# Scroll the view so row 28 is at the top and select A42, matching the
# author's final cursor/scroll position after adding the new rows.
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 1
$ws.Range("A42").Select()
